# Update countries & provincias Spain
# - Refresh the "datos actualizados" timestamp in A1.
# - Re-sort several rows of the country table (columns A..H) to reflect the
#   latest case counts: a handful of countries swap places with their
#   neighbours, and some of their stats (Casos activos / Muertes / etc.) get
#   updated with the newer figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 05:52"

# --- Country name (column A) reorderings -------------------------------
$ws.Cells.Item(166, 1).Value = "Mozambique"
$ws.Cells.Item(167, 1).Value = "Siria"

$ws.Cells.Item(171, 1).Value = "Mongolia"
$ws.Cells.Item(172, 1).Value = "Republica del Chad"
$ws.Cells.Item(173, 1).Value = "Guam"

$ws.Cells.Item(184, 1).Value = "Fiyi"
$ws.Cells.Item(185, 1).Value = "Islas Virgenes de los Estados Unidos"

$ws.Cells.Item(210, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(211, 1).Value = "Burundi"
$ws.Cells.Item(213, 1).Value = "Santo Tome y Principe"

$ws.Cells.Item(215, 1).Value = "Yemen"
$ws.Cells.Item(216, 1).Value = "San Pedro y Miquelon"

# --- Updated statistics (columns B..H) ----------------------------------
# Row 166/167 (Mozambique / Siria)
$ws.Cells.Item(166, 4).Value = 8
$ws.Cells.Item(166, 8).Value = 0

$ws.Cells.Item(167, 4).Value = 5
$ws.Cells.Item(167, 8).Value = 3

# Row 171/172/173 (Mongolia / Republica del Chad / Guam)
$ws.Cells.Item(171, 3).Value = 1
$ws.Cells.Item(171, 4).Value = 7
$ws.Cells.Item(171, 5).Value = 26

$ws.Cells.Item(172, 2).Value = 33
$ws.Cells.Item(172, 4).Value = 8
$ws.Cells.Item(172, 5).Value = 25
$ws.Cells.Item(172, 8).Value = 0

$ws.Cells.Item(173, 4).Value = 0
$ws.Cells.Item(173, 5).Value = 31
$ws.Cells.Item(173, 8).Value = 1

# Row 210/211/213 (Islas Virgenes Britanicas / Burundi / Santo Tome y Principe)
$ws.Cells.Item(210, 3).Value = 1
$ws.Cells.Item(210, 4).Value = 3
$ws.Cells.Item(210, 5).Value = 1

$ws.Cells.Item(211, 2).Value = 5
$ws.Cells.Item(211, 4).Value = 4
$ws.Cells.Item(211, 5).Value = 0
$ws.Cells.Item(211, 8).Value = 1

$ws.Cells.Item(213, 4).Value = 0
$ws.Cells.Item(213, 5).Value = 4
$ws.Cells.Item(213, 8).Value = 0
